$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# The "versjon" column (AR) for each of the three data rows is bumped
# from "0.0.1" to "0.1.0".
$ws.Range("AR2").Value = "0.1.0"
$ws.Range("AR3").Value = "0.1.0"
$ws.Range("AR4").Value = "0.1.0"

# Update the view: scroll the visible window over to the right (so column
# AI is the left-most visible column) and move the active selection to
# AR4 (bottom-right of the used range).
$win = $excel.ActiveWindow
$win.ScrollColumn = 35
$win.ScrollRow = 1
$ws.Range("AR4").Select() | Out-Null
